$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" (E) and "Valor Mora" (F) figures for the two 2307 / 2306
# period rows were entered swapped. Correct it by exchanging the values
# between row 16 and row 17.

$periodo16 = $ws.Range("E16").Value2
$periodo17 = $ws.Range("E17").Value2
$valor16   = $ws.Range("F16").Value2
$valor17   = $ws.Range("F17").Value2

$ws.Range("E16").Value = $periodo17
$ws.Range("E17").Value = $periodo16

$ws.Range("F16").Value = $valor17
$ws.Range("F17").Value = $valor16
